$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report now covers 2024-08-06 instead of 2024-08-05, and the previous
# 14:00 row (Hour 14) has been dropped - the sheet now ends at Hour 13,
# shrinking the used range from A1:T15 to A1:T14.
$ws.Rows.Item(15).Delete()

# Updated hourly station-load figures for 2024-08-06.
# Columns: Hour, S1..S8, Total, eS1..eS8, eTotal
$data = @(
    @(1, 19192, 5435, 8520, 0, 5088, 0, 15274, 14868, 68377, 19193.9192, 5435.5435, 8520.852, 0, 5088.5088, 0, 15275.5274, 14869.4868, 68383.8377),
    @(2, 18619, 5008, 8260, 0, 4996, 0, 14694, 14177, 65754, 18620.8619, 5008.5008, 8260.826, 0, 4996.4996, 0, 14695.4694, 14178.4177, 65760.5754),
    @(3, 17969, 4808, 8031, 0, 4678, 0, 14159, 13807, 63452, 17970.7969, 4808.4808, 8031.8031, 0, 4678.4678, 0, 14160.4159, 13808.3807, 63458.3452),
    @(4, 17223, 4700, 7899, 0, 4771, 0, 13738, 13447, 61778, 17224.7223, 4700.47, 7899.7899, 0, 4771.4771, 0, 13739.3738, 13448.3447, 61784.1778),
    @(5, 17898, 4889, 8098, 0, 4917, 0, 13890, 14009, 63701, 17899.7898, 4889.4889, 8098.8098, 0, 4917.4917, 0, 13891.389, 14010.4009, 63707.37009999999),
    @(6, 18664, 5231, 8714, 0, 5123, 0, 14547, 14642, 66921, 18665.8664, 5231.5231, 8714.8714, 0, 5123.5123, 0, 14548.4547, 14643.4642, 66927.6921),
    @(7, 18237, 4552, 8627, 0, 5517, 0, 14033, 14593, 65559, 18238.8237, 4552.4552, 8627.8627, 0, 5517.5517, 0, 14034.4033, 14594.4593, 65565.55589999999),
    @(8, 21576, 4438, 9964, 0, 9197, 0, 15693, 16434, 77302, 21578.1576, 4438.4438, 9964.9964, 0, 9197.9197, 0, 15694.5693, 16435.6434, 77309.7302),
    @(9, 26902, 5153, 12410, 0, 12608, 0, 18857, 18423, 94353, 26904.6902, 5153.5153, 12411.241, 0, 12609.2608, 0, 18858.8857, 18424.8423, 94362.4353),
    @(10, 27846, 5591, 17153, 0, 13638, 0, 0, 20468, 84696, 27848.7846, 5591.5591, 17154.7153, 0, 13639.3638, 0, 0, 20470.0468, 84704.4696),
    @(11, 28721, 5830, 18641, 0, 13864, 0, 17260, 21035, 105351, 28723.8721, 5830.583, 18642.8641, 0, 13865.3864, 0, 17261.726, 21037.1035, 105361.5351),
    @(12, 28957, 5853, 18818, 0, 14028, 0, 17640, 21320, 106616, 28959.8957, 5853.5853, 18819.8818, 0, 14029.4028, 0, 17641.764, 21322.132, 106626.6616),
    @(13, 27430, 5862, 19220, 0, 14041, 0, 17840, 21572, 105965, 27432.743, 5862.5862, 19221.922, 0, 14042.4041, 0, 17841.784, 21574.1572, 105975.5965)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2

    # Temporarily format column A as text so the "yyyy-mm-dd" date string is
    # stored literally instead of being auto-converted to a date serial,
    # then restore the cell's original (default) formatting.
    $dateCell = $ws.Range("A$r")
    $dateCell.NumberFormat = "@"
    $dateCell.Value2 = "2024-08-06"
    $dateCell.ClearFormats()

    $values = $data[$i]
    $arr = New-Object 'object[,]' 1,($values.Length)
    for ($j = 0; $j -lt $values.Length; $j++) { $arr[0,$j] = $values[$j] }
    $ws.Range("B${r}:T${r}").Value2 = $arr
}
